$wb = $excel.ActiveWorkbook

# --- Add "Truck_Amandla" sheet, positioned right after "Bus_Makhulu" -------
# Clone "Trailer_Elula" (same template layout/formatting) and drop it in
# place right after Bus_Makhulu, ahead of the trailer sheets.
$trailerElula = $wb.Worksheets.Item("Trailer_Elula")
$busMakhulu   = $wb.Worksheets.Item("Bus_Makhulu")
$insertAt     = $busMakhulu.Index + 1
$trailerElula.Copy([System.Reflection.Missing]::Value, $busMakhulu)
$truck = $wb.Worksheets.Item($insertAt)
$truck.Name = "Truck_Amandla"

$truck.Range("H3").Value = "Truck_Amandla"
$truck.Range("H4").Value = "sedan"
$truck.Range("H6").Value = 0.43
$truck.Range("F9").Value = -1.2
$truck.Range("G9").Value = 0
$truck.Range("H9").Value = 1.1000000000000001
[void]$truck.Range("H5:H9").Select()

# --- Add "Trailer_Kumanzi" sheet, positioned at the very end ---------------
$trailerThwala = $wb.Worksheets.Item("Trailer_Thwala")
$trailerThwala.Copy([System.Reflection.Missing]::Value, $trailerThwala)
$kumanzi = $wb.Worksheets.Item("Trailer_Thwala (2)")
$kumanzi.Name = "Trailer_Kumanzi"

$kumanzi.Range("H3").Value = "Trailer_Kumanzi"
$kumanzi.Range("H4").Value = "sedan"
$kumanzi.Range("H6").Value = 0.43
$kumanzi.Range("F9").Value = 5
$kumanzi.Range("G9").Value = 0
$kumanzi.Range("H9").Value = 2
[void]$kumanzi.Range("J20").Select()
